$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.809441918582433
$ws.Range("C2").Value = 1.773799324762276
$ws.Range("D2").Value = 6.212147121844095
$ws.Range("E2").Value = 1.700078601593715
$ws.Range("F2").Value = 1.807629892617897
$ws.Range("G2").Value = 2.002151662949331
$ws.Range("H2").Value = 1.700317937929734
$ws.Range("B3").Value = 1.820836472639407
$ws.Range("C3").Value = 1.790830877570518
$ws.Range("D3").Value = 4.809282045285278
$ws.Range("E3").Value = 1.703624045046675
$ws.Range("F3").Value = 1.817304492547184
$ws.Range("G3").Value = 1.971045113597998
$ws.Range("H3").Value = 1.703193982803089
$ws.Range("B4").Value = 1.782668010042132
$ws.Range("C4").Value = 1.860444782376231
$ws.Range("D4").Value = 5.544570533031403
$ws.Range("E4").Value = 1.689118343792704
$ws.Range("F4").Value = 1.782010133492829
$ws.Range("G4").Value = 1.993043150492126
$ws.Range("H4").Value = 1.689607354824597
$ws.Range("B5").Value = 1.817028306682256
$ws.Range("C5").Value = 1.998626836058858
$ws.Range("D5").Value = 2.934576195968682
$ws.Range("E5").Value = 1.705623931229836
$ws.Range("F5").Value = 1.812988199412726
$ws.Range("G5").Value = 2.020387323550141
$ws.Range("H5").Value = 1.705018684914926
$ws.Range("B6").Value = 1.823976775771101
$ws.Range("C6").Value = 1.996492174193518
$ws.Range("D6").Value = 0.8753565019908757
$ws.Range("E6").Value = 1.709811521561945
$ws.Range("F6").Value = 1.819960065524298
$ws.Range("G6").Value = 1.875463510908506
$ws.Range("H6").Value = 1.709059898859096
$ws.Range("B7").Value = 1.804522491512453
$ws.Range("C7").Value = 2.082638963899223
$ws.Range("D7").Value = 1.219932049956222
$ws.Range("E7").Value = 1.705657750267797
$ws.Range("F7").Value = 1.800608292030371
$ws.Range("G7").Value = 1.961212458634178
$ws.Range("H7").Value = 1.705004603960501
$ws.Range("B8").Value = 1.815118152670409
$ws.Range("C8").Value = 1.331268593591207
$ws.Range("D8").Value = 1.962454242750496
$ws.Range("E8").Value = 1.691569675101092
$ws.Range("F8").Value = 1.811121680113447
$ws.Range("G8").Value = 1.497617726625957
$ws.Range("H8").Value = 1.691113062267858
$ws.Range("B9").Value = 1.820110013471357
$ws.Range("C9").Value = 1.862728241502911
$ws.Range("D9").Value = 1.489277673447328
$ws.Range("E9").Value = 1.699117520271848
$ws.Range("F9").Value = 1.816956055849428
$ws.Range("G9").Value = 1.812437534362181
$ws.Range("H9").Value = 1.699213878649833
$ws.Range("B10").Value = 1.57578486270646
$ws.Range("C10").Value = 2.084409928598445
$ws.Range("D10").Value = 3.052325732473214
$ws.Range("E10").Value = 1.664463245616165
$ws.Range("F10").Value = 1.584180934600336
$ws.Range("G10").Value = 2.033644326230511
$ws.Range("H10").Value = 1.667191889118646
$ws.Range("B11").Value = 1.524747233310899
$ws.Range("C11").Value = 2.106358068430188
$ws.Range("D11").Value = 2.462184757981162
$ws.Range("E11").Value = 1.670736457539431
$ws.Range("F11").Value = 1.532756793977353
$ws.Range("G11").Value = 2.025370691254477
$ws.Range("H11").Value = 1.671843872822301
$ws.Range("B12").Value = 1.275774180011838
$ws.Range("C12").Value = 2.058589535148215
$ws.Range("D12").Value = 0.6115296961620851
$ws.Range("E12").Value = 1.617529040115271
$ws.Range("F12").Value = 1.280968368999795
$ws.Range("G12").Value = 1.892787997103477
$ws.Range("H12").Value = 1.617368168486356
$ws.Range("B13").Value = 1.559307860733568
$ws.Range("C13").Value = 2.079739547593903
$ws.Range("D13").Value = 2.167242324121294
$ws.Range("E13").Value = 1.658285052561531
$ws.Range("F13").Value = 1.565073904619898
$ws.Range("G13").Value = 1.989697056952977
$ws.Range("H13").Value = 1.659963471116276
